# Refresh crypto price/volume figures (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.250.51'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.594.16'
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.71'
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.94'
$ws.Range("E10").Value = '  -1.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0852'
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.818.57'
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.573.86'
$ws.Range("E13").Value = '  -0.64%  '
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("E15").Value = '  -2.38%  '
$ws.Range("E16").Value = '  -0.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.218.25'
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.15'
$ws.Range("E18").Value = '  +7.16%  '
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("E20").Value = '  +4.60%  '
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.90'
$ws.Range("E23").Value = '  -0.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.14'
$ws.Range("E24").Value = '  +1.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.62'
$ws.Range("E25").Value = '  +0.95%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("E28").Value = '  +0.68%  '
$ws.Range("E29").Value = '  +1.63%  '
$ws.Range("E30").Value = '  -0.38%  '
$ws.Range("E31").Value = '  +0.28%  '
$ws.Range("E32").Value = '  +0.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.461.47'
$ws.Range("E33").Value = '  +3.97%  '
$ws.Range("E34").Value = '  +0.46%  '
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("E37").Value = '  -4.08%  '
$ws.Range("E38").Value = '  -1.14%  '
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.79'
$ws.Range("E40").Value = '  -1.03%  '
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("E42").Value = '  +1.84%  '
$ws.Range("E43").Value = '  -1.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.730.95'
$ws.Range("E44").Value = '  +0.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.755'
$ws.Range("E45").Value = '  -1.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.41'
$ws.Range("E46").Value = '  -0.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '87.50'
$ws.Range("E47").Value = '  +2.63%  '
$ws.Range("E48").Value = '  -0.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.997'
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("E51").Value = '  -2.29%  '
